$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (row 1): Wins / Losses / Ties in columns AD, AE, AF
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header style (bold, centered, bordered) from the adjacent
# existing header cell (AC1) onto the three new header cells.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill in the team record (Wins=74, Losses=88, Ties=0) for every player
# row, 2 through 42.
for ($r = 2; $r -le 42; $r++) {
    $ws.Cells.Item($r, 30).Value = 74
    $ws.Cells.Item($r, 31).Value = 88
    $ws.Cells.Item($r, 32).Value = 0
}
